$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") on every existing data row (2..218) moves
#    from 2023-10-06 (45205) to 2023-10-07 (45206).
$ws.Range("C2:C218").Value = 45206

# 2) The former last row (218) now gets an explicit row height (matches
#    the new last row 219 that Excel stamps when a row is appended).
$ws.Rows.Item(218).RowHeight = 15

# 3) Append the new record as row 219.
$ws.Range("A219").Value = "A 48107-2023"
$ws.Range("B219").Value = 45205
$ws.Range("C219").Value = 45206
$ws.Range("B219:C219").NumberFormat = "YYYY-MM-DD"
$ws.Range("D219").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E219").Value = "ÖDESHÖG"
$ws.Range("G219").Value = 2.4
$ws.Range("H219").Value = 0
$ws.Range("I219").Value = 0
$ws.Range("J219").Value = 0
$ws.Range("K219").Value = 0
$ws.Range("L219").Value = 0
$ws.Range("M219").Value = 0
$ws.Range("N219").Value = 0
$ws.Range("O219").Value = 0
$ws.Range("P219").Value = 0
$ws.Range("Q219").Value = 0
$ws.Range("R219").WrapText = $true
